# Build site at 2023-04-12 14:53:07 UTC
# LOM3083.xlsx: fill in the real course-description texts (Objetivos, Programa
# resumido, Programa, Metodo, Criterio, Norma de recuperacao, Bibliografia),
# which previously held stray/misplaced values left over from the template,
# and add a dedicated row for 'Docentes responsaveis'.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivosText = 'Introdução de conceitos relacionados com taxa e fluxo de quantidade de movimento, calor e massa, aplicados ao processamento de materiais metálicos, cerâmicos e poliméricos. Capacitar o aluno a modelar e resolver problemas de interesse em fenômenos de transporte, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.'
$docenteText = '1176388 - Luiz Tadeu Fernandes Eleno'
$programaResumidoText = 'Introdução à transferência de calor. Condução de calor em regime permanente e em regime transiente. Transferência de calor por convecção livre e forçada. Transferência de calor por radiação térmica. Transferência de calor com mudança de fase. Transferência de massa. Exemplos de aplicação.'
$programaText = 'Introdução à transferência de calor. Propriedades térmicas dos materiais. Condutividade térmica de sólidos, fluidos e meios porosos. Conceito de difusividade térmica. 
Transferência de calor por condução: transferência de calor em regime permanente. Equação de Fourier. Transferência de calor em regime permanente com contornos convectivos. Lei de Newton do resfriamento. Condução de calor em regime transiente. Difusividade térmica. Número de Biot. Analogia entre transferência de calor e circuitos elétricos: conceitos de resistência e capacitância térmicas.
Transferência de calor por convecção livre e forçada. Convecção livre. Parâmetros de similaridade. Número de Rayleigh. Convecção forçada. Teoria da camada limite. Número de Prandtl e número de Nusselt.
Transferência de calor por radiação. Radiação do corpo negro. Propriedades da radiação. Fator de forma da radiação.
Transferência de calor na solidificação. 
Transferência de massa. Difusividade em sólidos, líquidos, gasosos e meios porosos. Transferência de massa em sistemas fluídos. Modelos para o coeficiente de transferência de massa. 
Transferência de calor com mudança de fase: ebulição e condensação.
Transferência de massa em sistemas heterogêneos. Reações sólidos/gás, sólido/líquido, líquido/líquido e líquido/gás.'
$metodoText = 'Aulas expositivas, seminários e exercícios comentados.'
$criterioText = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$normaText = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$biblioText = 'INCROPERA, F, P; DEWITT, D. P. Fundamentos de Transferência de Calor e de Massa, LTC Editora, 2005.
BENNETT, C. D.; MYERS, J. E. Fenômenos de Transporte. McGraw-Hill.
KREITH, F.; BOHN, M. S. Princípios de Transferência de Calor, Thomson Learning, 2003.
HOLMAN, J. P. Transferência de Calor, McGraw-Hill, 1983.
POIRIER, D.R.; GEIGER, G.H. Transport Phenomena in Materials Processing, TMS, 1994.
GASKELL, David R. Introduction to Transport Phenomena in Materials Engineering. Prentice Hall, 1991.
SZEKELY, J. Fluid Flow Phenomena in Metals Processing. Academic Press, 1979.'

# Insert a new row at 13 (old rows 13-23 shift down to 14-24, carrying their
# row heights along) to hold the 'Docentes responsaveis' value.
$ws.Rows(13).Insert()

# The insert leaves a stray formatted-but-empty A13 cell (column A's default
# style bleeding into the new row); the target layout has no A13 cell at all.
$ws.Cells.Item(13, 1).Clear()

# Objetivos (row 10): replace the misplaced professor name with the real text.
$ws.Cells.Item(10, 2).Value = $objetivosText
$ws.Cells.Item(10, 3).Value = $objetivosText

# Docentes responsaveis (new row 13): the professor name now lives here. Write
# the values first, then copy B20/C20's formatting (wrap text / red wrap text)
# onto the new cells so they reuse the existing styles instead of minting new ones.
$ws.Cells.Item(13, 2).Value = $docenteText
$ws.Cells.Item(13, 3).Value = $docenteText
$ws.Cells.Item(20, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(20, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Programa resumido (row 14): replace 'Semestral' with the short syllabus text.
$ws.Cells.Item(14, 2).Value = $programaResumidoText
$ws.Cells.Item(14, 3).Value = $programaResumidoText

# Programa (row 16): replace the stray date with the full syllabus text.
$ws.Cells.Item(16, 2).Value = $programaText
$ws.Cells.Item(16, 3).Value = $programaText

# Metodo (row 19): replace the misplaced professor name with the method text.
$ws.Cells.Item(19, 2).Value = $metodoText
$ws.Cells.Item(19, 3).Value = $metodoText

# Criterio (row 20): replace the method text with the grading-criteria text.
$ws.Cells.Item(20, 2).Value = $criterioText
$ws.Cells.Item(20, 3).Value = $criterioText

# Norma de recuperacao (row 21): replace the criteria text with the recovery text.
$ws.Cells.Item(21, 2).Value = $normaText
$ws.Cells.Item(21, 3).Value = $normaText

# Bibliografia (row 22): replace the recovery text with the actual bibliography.
$ws.Cells.Item(22, 2).Value = $biblioText
$ws.Cells.Item(22, 3).Value = $biblioText

# Column B no longer shares its width with column A; widen it to match column C,
# since it now holds the same long paragraph text.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

Write-Output "done"
